$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each row below mirrors one <row> changed in the refreshed cryptos export:
# price (D) and 1h-volume-change (E) refreshed; row 32/33 additionally swap
# which coin (PancakeSwap/Filecoin) occupies that rank.

$ws.Range("D2").Value = "34.355.49"
$ws.Range("E2").Value = "  -0.08%  "

$ws.Range("D3").Value = "1.801.08"
$ws.Range("E3").Value = "  +0.50%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").Value = "'227.40"
$ws.Range("E5").Value = "  +0.51%  "

$ws.Range("E6").Value = "  +3.91%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").Value = "'36.05"
$ws.Range("E8").Value = "  +10.31%  "

$ws.Range("E9").Value = "  +1.45%  "

$ws.Range("D10").Value = "'0.0693"
$ws.Range("E10").Value = "  +0.44%  "

$ws.Range("E11").Value = "  +2.05%  "

$ws.Range("D12").Value = "2.060.06"
$ws.Range("E12").Value = "  +0.52%  "

$ws.Range("E13").Value = "  +5.71%  "

$ws.Range("D14").Value = "1.784.83"
$ws.Range("E14").Value = "  -0.46%  "

$ws.Range("E15").Value = "  +1.30%  "

$ws.Range("D16").Value = "'4.50"
$ws.Range("E16").Value = "  +5.17%  "

$ws.Range("D17").Value = "34.346.04"
$ws.Range("E17").Value = "  -0.09%  "

$ws.Range("D18").Value = "'69.05"
$ws.Range("E18").Value = "  +1.01%  "

$ws.Range("D19").Value = "'245.20"
$ws.Range("E19").Value = "  +0.32%  "

$ws.Range("E20").Value = "  -0.24%  "

$ws.Range("D21").Value = "'11.66"
$ws.Range("E21").Value = "  +3.61%  "

$ws.Range("E22").Value = "  +0.07%  "

$ws.Range("D24").Value = "'2.15"
$ws.Range("E24").Value = "  +3.61%  "

$ws.Range("D25").Value = "'171.50"
$ws.Range("E25").Value = "  +3.01%  "

$ws.Range("D26").Value = "'7.97"
$ws.Range("E26").Value = "  +9.13%  "

$ws.Range("D27").Value = "'16.85"
$ws.Range("E27").Value = "  +2.13%  "

$ws.Range("E28").Value = "  +1.74%  "

$ws.Range("D30").Value = "'4.00"
$ws.Range("E30").Value = "  +0.46%  "

$ws.Range("E31").Value = "  +1.00%  "

$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "'1.24"
$ws.Range("E32").Value = "  +0.72%  "

$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'3.82"
$ws.Range("E33").Value = "  +0.27%  "

$ws.Range("D34").Value = "'1.81"
$ws.Range("E34").Value = "  +0.03%  "

$ws.Range("D35").Value = "1.398.20"
$ws.Range("E35").Value = "  -0.24%  "

$ws.Range("E36").Value = "  -0.81%  "

$ws.Range("D37").Value = "'2.46"
$ws.Range("E37").Value = "  -4.60%  "

$ws.Range("E38").Value = "  -0.08%  "

$ws.Range("E39").Value = "  -0.33%  "

$ws.Range("D40").Value = "'1.23"
$ws.Range("E40").Value = "  +10.44%  "

$ws.Range("E41").Value = "  +2.70%  "

$ws.Range("D42").Value = "'82.58"
$ws.Range("E42").Value = "  -2.88%  "

$ws.Range("D43").Value = "'2.83"
$ws.Range("E43").Value = "  -0.28%  "

$ws.Range("E44").Value = "  +0.49%  "

$ws.Range("D45").Value = "'13.41"
$ws.Range("E45").Value = "  -3.03%  "

$ws.Range("D46").Value = "'0.0508"
$ws.Range("E46").Value = "  -3.45%  "

$ws.Range("E47").Value = "  +0.07%  "

$ws.Range("D48").Value = "1.961.16"
$ws.Range("E48").Value = "  +0.70%  "

$ws.Range("D49").Value = "'104.27"
$ws.Range("E49").Value = "  -0.81%  "

$ws.Range("E50").Value = "  +0.09%  "

$ws.Range("D51").Value = "0.0₆0128"
$ws.Range("E51").Value = "  +0.23%  "
